$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the run we need to split:
#   " pour recepvoir en commu"
# which must become:
#   " pour " + "<del>" + "avoyr en" + "</del>" + " " + "recepvoir en commu"
# (the first / last pieces keep the original run's formatting; the
#  middle pieces get new formatting matching markup already used
#  elsewhere in the document: a red Courier-New style for the
#  "<del>"/"</del>" tags, and a bare/plain style for the deleted
#  text and the following space)
# ------------------------------------------------------------------

$target = $d.Content
$found = $target.Find.Execute(" pour recepvoir en commu", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "could not find target text"
}

# split point sits right after ' pour ' (6 chars) and before 'recepvoir...'
$splitPoint = $target.Start + 6

# ------------------------------------------------------------------
# Piece 1: "<del>" with markup formatting (rFonts Courier New x4,
# color a91111, sz/szCs 18). We copy an existing "<corr>" run (same
# formatting is used for that markup elsewhere in the document) and
# paste it in, then overwrite its text -- this is the only way to
# get a full <w:rFonts .../> (all four attributes) out of this
# COM-interop host, since Font.Name only ever emits ascii+hAnsi.
# ------------------------------------------------------------------
$markupSrc1 = $d.Content
$markupFound1 = $markupSrc1.Find.Execute("<corr>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $markupFound1) {
    throw "could not find markup source run (1)"
}
$markupSrc1.Copy()

$insPoint1 = $d.Range($splitPoint, $splitPoint)
$insPoint1.Paste()
$pastedRange1 = $d.Range($splitPoint, $splitPoint + 6)
$pastedRange1.Text = "<del>"
$afterPiece1 = $splitPoint + 5

# ------------------------------------------------------------------
# Piece 2: "avoyr en" with plain formatting (just rtl, no color or
# font override). Copy source: the single-letter run rendering
# "left-top" elsewhere in the doc, which already has exactly that
# bare rPr.
# ------------------------------------------------------------------
$plainSrc1 = $d.Content
$plainFound1 = $plainSrc1.Find.Execute("left-top", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $plainFound1) {
    throw "could not find plain source run (1)"
}
$plainSrc1.Copy()

$insPoint2 = $d.Range($afterPiece1, $afterPiece1)
$insPoint2.Paste()
$pastedRange2 = $d.Range($afterPiece1, $afterPiece1 + 8)
$pastedRange2.Text = "avoyr en"
$afterPiece2 = $afterPiece1 + 8

# ------------------------------------------------------------------
# Piece 3: "</del>" with markup formatting -- same recipe as piece 1.
# ------------------------------------------------------------------
$markupSrc2 = $d.Content
$markupFound2 = $markupSrc2.Find.Execute("<corr>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $markupFound2) {
    throw "could not find markup source run (2)"
}
$markupSrc2.Copy()

$insPoint3 = $d.Range($afterPiece2, $afterPiece2)
$insPoint3.Paste()
$pastedRange3 = $d.Range($afterPiece2, $afterPiece2 + 6)
$pastedRange3.Text = "</del>"
$afterPiece3 = $afterPiece2 + 6

# ------------------------------------------------------------------
# Piece 4: " " with plain formatting -- same recipe as piece 2.
# ------------------------------------------------------------------
$plainSrc2 = $d.Content
$plainFound2 = $plainSrc2.Find.Execute("left-top", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $plainFound2) {
    throw "could not find plain source run (2)"
}
$plainSrc2.Copy()

$insPoint4 = $d.Range($afterPiece3, $afterPiece3)
$insPoint4.Paste()
$pastedRange4 = $d.Range($afterPiece3, $afterPiece3 + 8)
$pastedRange4.Text = " "
$afterPiece4 = $afterPiece3 + 1

Write-Host "Edit applied successfully."
